$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.820.89'
$ws.Range("E2").Value = '  +1.19%  '

$ws.Range("D3").Value = '2.085.95'
$ws.Range("E3").Value = '  +0.99%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.45'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.623'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.40%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.27'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.88%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.387'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.51%  '

$ws.Range("E10").Value = '  +2.23%  '

$ws.Range("E11").Value = '  +2.90%  '

$ws.Range("D12").Value = '2.383.89'
$ws.Range("E12").Value = '  +0.67%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.39'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.03'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.73%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.761'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.10%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.23'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.24%  '

$ws.Range("D17").Value = '2.088.34'
$ws.Range("E17").Value = '  +1.17%  '

$ws.Range("D18").Value = '37.697.75'
$ws.Range("E18").Value = '  +1.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.13'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.38%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.86'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.03%  '

$ws.Range("D21").Value = '0.0₃0820'
$ws.Range("E21").Value = '  +1.36%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.82'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.69%  '

$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("E24").Value = '  -1.85%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.69%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.02'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.139'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +9.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.91'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.49%  '

$ws.Range("E29").Value = '  -0.39%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.45'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.05%  '

$ws.Range("E31").Value = '  +0.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.60'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0624'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.31%  '

$ws.Range("E34").Value = '  +0.22%  '

$ws.Range("E35").Value = '  +0.65%  '

$ws.Range("E36").Value = '  +3.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.37'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.19%  '

$ws.Range("E38").Value = '  +0.15%  '

$ws.Range("E39").Value = '  -4.74%  '

$ws.Range("E40").Value = '  +6.36%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.93'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.78%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.99'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.74%  '

$ws.Range("E43").Value = '  +0.43%  '

$ws.Range("D44").Value = '1.452.13'
$ws.Range("E44").Value = '  -0.85%  '

$ws.Range("E45").Value = '  -0.86%  '

$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.05'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.30%  '

$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.08'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -5.86%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.59'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.64%  '

$ws.Range("E49").Value = '  +2.88%  '

$ws.Range("E50").Value = '  +1.73%  '

$ws.Range("D51").Value = '2.278.95'
$ws.Range("E51").Value = '  +1.10%  '
